$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.686.59"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "'1.599.97"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'211.10"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.0617"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "'19.62"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "'1.824.55"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "'1.607.42"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "'64.79"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "'26.670.82"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'208.87"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'6.78"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'2.32"
$ws.Range("E23").Value = "  -3.32%  "
$ws.Range("D24").Value = "'8.90"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'145.72"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").Value = "'15.26"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "'0.661"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "'1.296.27"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "'0.845"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'5.40"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "'63.71"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "'1.736.99"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").Value = "'0.898"
$ws.Range("E46").Value = "  +7.29%  "
$ws.Range("D47").Value = "'90.14"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "'0.0505"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = "  +0.20%  "
